$wb = $excel.ActiveWorkbook

# Map of row -> [F new value, G new value (optional)]
$updates = @{
    3  = @{ F = 3060 }
    5  = @{ F = 158 }
    7  = @{ F = 1681 }
    11 = @{ F = 3 }
    12 = @{ F = 1377 }
    14 = @{ F = 523 }
    16 = @{ F = 35 }
    17 = @{ F = 6 }
    21 = @{ F = 91 }
    22 = @{ F = 106; G = 55 }
    23 = @{ F = 3224 }
    24 = @{ F = 393 }
    25 = @{ F = 142 }
    26 = @{ F = 324 }
    27 = @{ F = 10 }
    29 = @{ F = 97 }
}

$sheetNames = @("展览", "全部类型")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)
    foreach ($row in $updates.Keys) {
        $vals = $updates[$row]
        $ws.Range("F$row").Value = $vals.F
        if ($vals.ContainsKey("G")) {
            $ws.Range("G$row").Value = $vals.G
        }
    }
}
